$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Remove the row that becomes fully absorbed into row 12 after the merge
$t.Rows.Item(13).Delete()

# Append a new row for the added math-fact row at the end of the table
$t.Rows.Add() | Out-Null

# Write every cell value for the final 20x5 grid
$t.Cell(1, 1).Range.Text = "91-22="
$t.Cell(1, 2).Range.Text = "6+26="
$t.Cell(1, 3).Range.Text = "90-22="
$t.Cell(1, 4).Range.Text = "38+54="
$t.Cell(1, 5).Range.Text = "42-5="
$t.Cell(2, 1).Range.Text = "80-46="
$t.Cell(2, 2).Range.Text = "84-76="
$t.Cell(2, 3).Range.Text = "40-32="
$t.Cell(2, 4).Range.Text = "4+59="
$t.Cell(2, 5).Range.Text = "95-69="
$t.Cell(3, 1).Range.Text = "16+77="
$t.Cell(3, 2).Range.Text = "71-16="
$t.Cell(3, 3).Range.Text = "22-13="
$t.Cell(3, 4).Range.Text = "19+72="
$t.Cell(3, 5).Range.Text = "46+27="
$t.Cell(4, 1).Range.Text = "27+47="
$t.Cell(4, 2).Range.Text = "19+29="
$t.Cell(4, 3).Range.Text = "37+27="
$t.Cell(4, 4).Range.Text = "25-9="
$t.Cell(4, 5).Range.Text = "37+37="
$t.Cell(5, 1).Range.Text = "8+24="
$t.Cell(5, 2).Range.Text = "26+37="
$t.Cell(5, 3).Range.Text = "87-78="
$t.Cell(5, 4).Range.Text = "37+5="
$t.Cell(5, 5).Range.Text = "63+8="
$t.Cell(6, 1).Range.Text = "27+48="
$t.Cell(6, 2).Range.Text = "48+3="
$t.Cell(6, 3).Range.Text = "51-6="
$t.Cell(6, 4).Range.Text = "88-9="
$t.Cell(6, 5).Range.Text = "97-19="
$t.Cell(7, 1).Range.Text = "21-6="
$t.Cell(7, 2).Range.Text = "19+23="
$t.Cell(7, 3).Range.Text = "61-16="
$t.Cell(7, 4).Range.Text = "10-6="
$t.Cell(7, 5).Range.Text = "8+43="
$t.Cell(8, 1).Range.Text = "18+63="
$t.Cell(8, 2).Range.Text = "98-49="
$t.Cell(8, 3).Range.Text = "39+23="
$t.Cell(8, 4).Range.Text = "33+28="
$t.Cell(8, 5).Range.Text = "65-18="
$t.Cell(9, 1).Range.Text = "42+39="
$t.Cell(9, 2).Range.Text = "40-37="
$t.Cell(9, 3).Range.Text = "44-38="
$t.Cell(9, 4).Range.Text = "10-9="
$t.Cell(9, 5).Range.Text = "82-14="
$t.Cell(10, 1).Range.Text = "7+68="
$t.Cell(10, 2).Range.Text = "66+15="
$t.Cell(10, 3).Range.Text = "13+19="
$t.Cell(10, 4).Range.Text = "23-14="
$t.Cell(10, 5).Range.Text = "28+24="
$t.Cell(11, 1).Range.Text = "72-63="
$t.Cell(11, 2).Range.Text = "58+15="
$t.Cell(11, 3).Range.Text = "48+27="
$t.Cell(11, 4).Range.Text = "83-56="
$t.Cell(11, 5).Range.Text = "91-69="
$t.Cell(12, 1).Range.Text = "56+39="
$t.Cell(12, 2).Range.Text = "37+6="
$t.Cell(12, 3).Range.Text = "96-49="
$t.Cell(12, 4).Range.Text = "53+39="
$t.Cell(12, 5).Range.Text = "72-5="
$t.Cell(13, 1).Range.Text = "72-53="
$t.Cell(13, 2).Range.Text = "32-15="
$t.Cell(13, 3).Range.Text = "68+18="
$t.Cell(13, 4).Range.Text = "27+6="
$t.Cell(13, 5).Range.Text = "72-17="
$t.Cell(14, 1).Range.Text = "54+29="
$t.Cell(14, 2).Range.Text = "8+29="
$t.Cell(14, 3).Range.Text = "21-9="
$t.Cell(14, 4).Range.Text = "81-24="
$t.Cell(14, 5).Range.Text = "19+76="
$t.Cell(15, 1).Range.Text = "92-34="
$t.Cell(15, 2).Range.Text = "94-59="
$t.Cell(15, 3).Range.Text = "48+39="
$t.Cell(15, 4).Range.Text = "59+37="
$t.Cell(15, 5).Range.Text = "90-17="
$t.Cell(16, 1).Range.Text = "21-17="
$t.Cell(16, 2).Range.Text = "27+69="
$t.Cell(16, 3).Range.Text = "59+8="
$t.Cell(16, 4).Range.Text = "9+2="
$t.Cell(16, 5).Range.Text = "79+18="
$t.Cell(17, 1).Range.Text = "49+4="
$t.Cell(17, 2).Range.Text = "36+46="
$t.Cell(17, 3).Range.Text = "45+39="
$t.Cell(17, 4).Range.Text = "5+7="
$t.Cell(17, 5).Range.Text = "22-5="
$t.Cell(18, 1).Range.Text = "50-14="
$t.Cell(18, 2).Range.Text = "60-37="
$t.Cell(18, 3).Range.Text = "70-47="
$t.Cell(18, 4).Range.Text = "83-64="
$t.Cell(18, 5).Range.Text = "71-2="
$t.Cell(19, 1).Range.Text = "75-29="
$t.Cell(19, 2).Range.Text = "70-62="
$t.Cell(19, 3).Range.Text = "62-27="
$t.Cell(19, 4).Range.Text = "27+8="
$t.Cell(19, 5).Range.Text = "56-27="
$t.Cell(20, 1).Range.Text = "12-8="
$t.Cell(20, 2).Range.Text = "61-14="
$t.Cell(20, 3).Range.Text = "58+25="
$t.Cell(20, 4).Range.Text = "71-66="
$t.Cell(20, 5).Range.Text = "32-24="

Write-Output "done"
